# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.047.66'
$ws.Range("E2").Value = '  +6.64%  '
$ws.Range("D3").Value = '3.017.56'
$ws.Range("E3").Value = '  +4.01%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'584.04"
$ws.Range("E5").Value = '  +2.59%  '
$ws.Range("D6").Value = "'163.46"
$ws.Range("E6").Value = '  +13.51%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '3.012.65'
$ws.Range("E8").Value = '  +3.92%  '
$ws.Range("E9").Value = '  +3.27%  '
$ws.Range("E10").Value = '  +1.34%  '
$ws.Range("E11").Value = '  +7.90%  '
$ws.Range("E12").Value = '  +6.43%  '
$ws.Range("E13").Value = '  +9.07%  '
$ws.Range("D14").Value = "'34.96"
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '66.064.21'
$ws.Range("E16").Value = '  +6.79%  '
$ws.Range("D17").Value = '3.519.02'
$ws.Range("E17").Value = '  +4.07%  '
$ws.Range("D18").Value = "'6.99"
$ws.Range("E18").Value = '  +7.15%  '
$ws.Range("D19").Value = '3.016.92'
$ws.Range("E19").Value = '  +4.06%  '
$ws.Range("D20").Value = "'458.40"
$ws.Range("E20").Value = '  +6.11%  '
$ws.Range("D21").Value = "'13.99"
$ws.Range("E21").Value = '  +7.97%  '
$ws.Range("E22").Value = '  +5.42%  '
$ws.Range("E23").Value = '  +7.70%  '
$ws.Range("D24").Value = "'82.51"
$ws.Range("E24").Value = '  +4.61%  '
$ws.Range("E26").Value = '  +2.65%  '
$ws.Range("D27").Value = "'10.64"
$ws.Range("E27").Value = '  +5.12%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = "'8.19"
$ws.Range("E29").Value = '  +16.72%  '
$ws.Range("E30").Value = '  +15.02%  '
$ws.Range("E31").Value = '  +4.37%  '
$ws.Range("E32").Value = '  -6.58%  '
$ws.Range("E33").Value = '  +5.66%  '
$ws.Range("E34").Value = '  +3.47%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = "'0.995"
$ws.Range("E36").Value = '  +4.10%  '
$ws.Range("E37").Value = '  +7.63%  '
$ws.Range("E38").Value = '  +11.61%  '
$ws.Range("D39").Value = "'3.05"
$ws.Range("E39").Value = '  +6.97%  '
$ws.Range("D40").Value = "'49.95"
$ws.Range("E40").Value = '  +2.26%  '
$ws.Range("D41").Value = "'0.310"
$ws.Range("E41").Value = '  +15.06%  '
$ws.Range("E42").Value = '  +6.28%  '
$ws.Range("D43").Value = "'43.71"
$ws.Range("E43").Value = '  +8.51%  '
$ws.Range("E44").Value = '  +4.28%  '
$ws.Range("D45").Value = "'389.55"
$ws.Range("E45").Value = '  +12.36%  '
$ws.Range("E46").Value = '  +6.60%  '
$ws.Range("D47").Value = '2.798.61'
$ws.Range("E47").Value = '  +3.69%  '
$ws.Range("D48").Value = "'135.41"
$ws.Range("E48").Value = '  +2.78%  '
$ws.Range("D50").Value = "'24.10"
$ws.Range("E50").Value = '  +11.48%  '
$ws.Range("E51").Value = '  +4.16%  '
